# "regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals" -- the pitching log's strikeout column (column G, header
# "K") is regenerated with freshly (re)computed per-appearance strikeout
# counts. Column G is index 7 (A=1 ... J=10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4;  3  = 0;  4  = 1;  5  = 2;  6  = 5;  7  = 1;  8  = 2;  9  = 3;
    10 = 1;  11 = 1;  12 = 1;  13 = 3;  14 = 1;  15 = 2;  16 = 1;  17 = 1;
    18 = 1;  19 = 1;  20 = 0;  21 = 1;  22 = 0;  23 = 1;  25 = 0;  26 = 4;
    27 = 1;  28 = 1;  29 = 1;  30 = 3;  32 = 1;  33 = 1;  34 = 2;  35 = 4;
    36 = 1;  37 = 1;  38 = 2;  39 = 1;  40 = 2;  41 = 1;  42 = 2;  44 = 1;
    45 = 0;  47 = 1;  48 = 1;  49 = 0;  50 = 2;  51 = 1;  52 = 0;  53 = 0;
    54 = 1;  55 = 1;  56 = 1;  57 = 1;  58 = 2;  59 = 0;  60 = 1;  61 = 1;
    62 = 1;  63 = 1;  64 = 2;  65 = 0;  66 = 1;  67 = 2;  68 = 1;  69 = 1;
    70 = 2;  71 = 2;  72 = 1;  73 = 2;  74 = 1;  75 = 0;  76 = 3;  77 = 3;
    78 = 3;  79 = 1;  80 = 0;  81 = 1;  82 = 3;  83 = 1;  84 = 1;  85 = 0;
    86 = 2;  87 = 2;  89 = 3;  90 = 1;  91 = 3;  92 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
